# Updates the "cryptos" worksheet with refreshed price / volume figures
# (and a re-ranking swap between Kaspa and Binance-PegBSC-USD at rows 35/36),
# matching a newer pull of the coinranking.com data.
#
# All of the Price (D) / Volume(1h) (E) / Coin (B) / Link (C) columns are
# stored as literal text in the sheet (e.g. "69.342.35", "  -2.60%  "), so
# writes are forced to text (leading apostrophe) and the style is reset back
# to Normal afterwards so Excel's "looks like a number" auto-conversion
# doesn't turn them into numeric/date cells or leave a stray quote-prefixed
# style behind.

function Set-CellText($sheet, $addr, $text) {
    $sheet.Range($addr).Value = "'" + $text
    $sheet.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "69.342.35"
Set-CellText $ws "E2" "  -2.60%  "

Set-CellText $ws "D3" "3.684.01"
Set-CellText $ws "E3" "  -3.40%  "

Set-CellText $ws "E4" "  +0.02%  "

Set-CellText $ws "D5" "688.59"
Set-CellText $ws "E5" "  -2.55%  "

Set-CellText $ws "D6" "162.50"
Set-CellText $ws "E6" "  -5.48%  "

Set-CellText $ws "D7" "3.683.24"
Set-CellText $ws "E7" "  -3.40%  "

Set-CellText $ws "E8" "  +0.07%  "

Set-CellText $ws "E9" "  -4.68%  "

Set-CellText $ws "E10" "  -8.52%  "

Set-CellText $ws "D11" "7.37"
Set-CellText $ws "E11" "  -4.14%  "

Set-CellText $ws "D12" "0.444"
Set-CellText $ws "E12" "  -3.85%  "

Set-CellText $ws "E13" "  -5.25%  "

Set-CellText $ws "D14" "33.53"
Set-CellText $ws "E14" "  -6.76%  "

Set-CellText $ws "D15" "4.305.52"
Set-CellText $ws "E15" "  -3.40%  "

Set-CellText $ws "D16" "3.687.49"
Set-CellText $ws "E16" "  -3.02%  "

Set-CellText $ws "D17" "69.407.00"
Set-CellText $ws "E17" "  -2.47%  "

Set-CellText $ws "D19" "16.32"
Set-CellText $ws "E19" "  -7.02%  "

Set-CellText $ws "D20" "6.62"
Set-CellText $ws "E20" "  -7.73%  "

Set-CellText $ws "D21" "482.45"
Set-CellText $ws "E21" "  -5.99%  "

Set-CellText $ws "D22" "9.93"
Set-CellText $ws "E22" "  -7.05%  "

Set-CellText $ws "E23" "  -7.68%  "

Set-CellText $ws "D24" "80.18"
Set-CellText $ws "E24" "  -4.76%  "

Set-CellText $ws "D25" "3.830.16"

Set-CellText $ws "E26" "  -9.48%  "

Set-CellText $ws "E27" "  +0.06%  "

Set-CellText $ws "E28" "  -4.90%  "

Set-CellText $ws "E29" "  -8.57%  "

Set-CellText $ws "E30" "  -10.15%  "

Set-CellText $ws "E31" "  -10.38%  "

Set-CellText $ws "D32" "6.85"
Set-CellText $ws "E32" "  -7.65%  "

Set-CellText $ws "E33" "  -7.47%  "

Set-CellText $ws "D34" "27.14"
Set-CellText $ws "E34" "  -6.81%  "

# Rows 35/36 swap ranking order: Binance-PegBSC-USD moves up to 35,
# Kaspa moves down to 36, each with refreshed price/volume data.
Set-CellText $ws "B35" "Binance-PegBSC-USD"
Set-CellText $ws "C35" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-CellText $ws "D35" "1.00"
Set-CellText $ws "E35" "  +0.09%  "

Set-CellText $ws "B36" "Kaspa"
Set-CellText $ws "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText $ws "D36" "0.165"
Set-CellText $ws "E36" "  -3.86%  "

Set-CellText $ws "D37" "3.650.87"
Set-CellText $ws "E37" "  -3.30%  "

Set-CellText $ws "D38" "8.52"
Set-CellText $ws "E38" "  -7.12%  "

Set-CellText $ws "D39" "6.31"
Set-CellText $ws "E39" "  +5.18%  "

Set-CellText $ws "D40" "2.33"
Set-CellText $ws "E40" "  -2.21%  "

Set-CellText $ws "D41" "0.0932"
Set-CellText $ws "E41" "  -7.71%  "

Set-CellText $ws "E43" "  -0.03%  "

Set-CellText $ws "E44" "  -7.09%  "

Set-CellText $ws "D45" "163.54"
Set-CellText $ws "E45" "  -3.90%  "

Set-CellText $ws "D46" "47.97"
Set-CellText $ws "E46" "  -3.31%  "

Set-CellText $ws "E47" "  -13.26%  "

Set-CellText $ws "D48" "29.85"
Set-CellText $ws "E48" "  +2.71%  "

Set-CellText $ws "E49" "  +0.70%  "

Set-CellText $ws "D50" "0.000288"
Set-CellText $ws "E50" "  -7.46%  "

Set-CellText $ws "E51" "  -0.63%  "
